# Daily price-tracker update: a new day's row is inserted at the top of the
# data (row 2), pushing all the existing date rows down by one. The new
# row gets the latest date with the same (unchanged) price figures that
# every other row already has. The sheet's used range therefore grows from
# A1:D78 to A1:D79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data rows (old row 2..78) down by one to make room for the
# newest date at row 2.
$ws.Rows.Item(2).Insert()

# Column A stores dates as plain text (not real Excel dates) throughout the
# sheet, so force the cell to Text before writing the literal date string -
# otherwise Excel auto-converts the "yyyy-mm-dd" text into a date serial.
# Clear the formatting right back afterwards so the cell ends up with the
# same (default/general) style as every other date cell in the column.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-06"
$ws.Range("A2").ClearFormats()

# The tracked commodity prices are unchanged day over day.
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
